$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename existing "Lüderitz" label to include the first demand tier suffix
$ws.Cells.Item(2, 1).Value = "Lüderitz 5400 kga"

# Update the existing row's annual demand value (column D) to the new first tier
$ws.Cells.Item(2, 4).Value = 5400

# Lat/lon values shared across the new demand-tier rows (same site, Lüderitz)
$lat = -26.642877645011101
$lon = 15.1439290700957

# Demand tiers to append as new rows 3-7 (name suffix, row, demand value)
$tiers = @(
    @{ Row = 3; Suffix = "54000"; Demand = 54000 },
    @{ Row = 4; Suffix = "540000"; Demand = 540000 },
    @{ Row = 5; Suffix = "5400000"; Demand = 5400000 },
    @{ Row = 6; Suffix = "54000000"; Demand = 54000000 },
    @{ Row = 7; Suffix = "540000000"; Demand = 540000000 }
)

foreach ($tier in $tiers) {
    $r = $tier.Row

    $ws.Cells.Item($r, 1).Value = "Lüderitz " + $tier.Suffix + " kga"

    $ws.Cells.Item($r, 2).Value = $lat
    $ws.Cells.Item($r, 2).NumberFormat = "0.00"

    $ws.Cells.Item($r, 3).Value = $lon
    $ws.Cells.Item($r, 3).NumberFormat = "0.00"

    $ws.Cells.Item($r, 4).Value = $tier.Demand
    $ws.Cells.Item($r, 4).NumberFormat = "_(* #,##0_);_(* \(#,##0\);_(* ""-""??_);_(@_)"

    $ws.Cells.Item($r, 5).Value = "NH3"
}

# Move the selection, matching the author's final cursor position
$ws.Range("A9").Select()
